$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B17").Value = 5
$ws.Range("B18").Select()
